# Apply the edits described by the commit "added data and updated images draft"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell content updates on the shelter data table (A3:F9) ---

# Row 4: Central Okanagan Foundation
$ws.Range("C4").Value = "8 am to 6 pm weekdays "
$ws.Range("F4").Value = "Food and homelessness "

# Row 5: "Cornerstore Shelter " -> "Cornerstone Shelter " (typo fix)
$ws.Range("A5").Value = "Cornerstone Shelter "

# Row 6: Kelowna's Gospel Mission
$ws.Range("C6").Value = "8 - 4 weekdays"
$ws.Range("F6").Value = "Social services organization, homelessness"

# Row 7: Harvey House - John Howard Society (time of opening was previously empty)
$ws.Range("C7").Value = "Open now(not mentioned)"

# Row 8: United Way British Columbia
$ws.Range("C8").Value = "8:30 - 4:30 weekdays"

# --- Update the active selection on the sheet to reflect the saved view state ---
$ws.Range("B14").Select()
